# Refresh the cryptos list snapshot (price + 1h volume change) with the
# latest values pulled from coinranking.com, as the GitHub Actions job
# does on every run. Rows 29/30 (PEPE / Binance-PegBSC-USD) also swapped
# rank order this time, so their Coin/Link/Price/Volume cells move together.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.315.33'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '2.617.96'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''596.11'  # keep as text, not a number
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = '''152.91'  # keep as text, not a number
$ws.Range('E6').Value = '  -1.63%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +2.63%  '
$ws.Range('D9').Value = '2.617.61'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  -2.82%  '
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').Value = '''5.18'  # keep as text, not a number
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('D14').Value = '''27.71'  # keep as text, not a number
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').Value = '3.096.72'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('E16').Value = '  -3.69%  '
$ws.Range('D17').Value = '67.282.86'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').Value = '2.614.44'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').Value = '''11.09'  # keep as text, not a number
$ws.Range('E19').Value = '  -2.11%  '
$ws.Range('D20').Value = '''363.12'  # keep as text, not a number
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('D21').Value = '''7.46'  # keep as text, not a number
$ws.Range('E21').Value = '  -4.37%  '
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('E23').Value = '  +1.58%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').Value = '''71.10'  # keep as text, not a number
$ws.Range('E25').Value = '  +4.74%  '
$ws.Range('D26').Value = '''10.00'  # keep as text, not a number
$ws.Range('E26').Value = '  -1.76%  '
$ws.Range('D27').Value = '2.743.53'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').Value = '''587.55'  # keep as text, not a number
$ws.Range('E28').Value = '  -3.83%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = '''1.00'  # keep as text, not a number
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '''0.0000102'  # keep as text, not a number
$ws.Range('E30').Value = '  -2.44%  '
$ws.Range('E31').Value = '  -3.99%  '
$ws.Range('D32').Value = '''7.83'  # keep as text, not a number
$ws.Range('E32').Value = '  -2.17%  '
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('E36').Value = '  -2.73%  '
$ws.Range('D37').Value = '''4.88'  # keep as text, not a number
$ws.Range('E37').Value = '  -1.95%  '
$ws.Range('D38').Value = '''157.46'  # keep as text, not a number
$ws.Range('D39').Value = '''19.10'  # keep as text, not a number
$ws.Range('E39').Value = '  -2.52%  '
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('D41').Value = '''5.26'  # keep as text, not a number
$ws.Range('E41').Value = '  -3.79%  '
$ws.Range('E42').Value = '  -2.11%  '
$ws.Range('D43').Value = '''2.56'  # keep as text, not a number
$ws.Range('E43').Value = '  -2.57%  '
$ws.Range('D44').Value = '''41.13'  # keep as text, not a number
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '''16.36'  # keep as text, not a number
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('D47').Value = '''156.81'  # keep as text, not a number
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('D48').Value = '0.0₆0289'
$ws.Range('E48').Value = '  -1.24%  '
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('D50').Value = '''21.91'  # keep as text, not a number
$ws.Range('E50').Value = '  +4.57%  '
$ws.Range('D51').Value = '''0.622'  # keep as text, not a number
$ws.Range('E51').Value = '  -0.57%  '
